# Applies the scheduled-runner value refresh to the Leve profit sheets.
# Each sheet has independent currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) that get re-pulled from the market-board source; this
# script pokes in the refreshed literal values cell-by-cell.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    ,("H17", 827.1429000000001)
    ,("J17", 818.5)
    ,("L17", 2455.5)
    ,("N17", -2791.5)
    ,("H31", 800)
    ,("I31", 0)
    ,("K31", 0)
    ,("H33", 245.87878)
    ,("J33", 507.55554)
    ,("L33", 507.55554)
    ,("N33", -965.5555400000001)
    ,("H88", 5774.1816)
    ,("I88", 1980)
    ,("J88", 8936)
    ,("K88", 1980)
    ,("L88", 8936)
    ,("M88", -1574)
    ,("N88", -9748)
    ,("H91", 5774.1816)
    ,("I91", 1980)
    ,("J91", 8936)
    ,("K91", 1980)
    ,("L91", 8936)
    ,("M91", -576)
    ,("N91", -11744)
    ,("H98", 1813.0769)
    ,("I98", 1088.9584)
    ,("K98", 1088.9584)
    ,("M98", 409.0416)
    ,("H122", 1813.0769)
    ,("I122", 1088.9584)
    ,("K122", 3266.8752)
    ,("M122", -816.8751999999999)
    ,("H137", 23134.74)
    ,("I137", 17551.055)
    ,("J137", 34302.11)
    ,("K137", 52653.165)
    ,("L137", 102906.33)
    ,("M137", -50103.165)
    ,("N137", -108006.33)
    ,("H138", 26901.65)
    ,("I138", 1587.8235)
    ,("K138", 4763.470499999999)
    ,("M138", 376.5295000000006)
)
foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
$clears = @("M31")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    ,("H2", 2696.4443)
    ,("I2", 2769.261)
    ,("J2", 2277.75)
    ,("K2", 2769.261)
    ,("L2", 2277.75)
    ,("M2", -2656.261)
    ,("N2", -2503.75)
    ,("H45", 3702.389)
    ,("I45", 2065)
    ,("K45", 2065)
    ,("M45", -1688)
    ,("H61", 7285.2354)
    ,("I61", 1297.4615)
    ,("K61", 1297.4615)
    ,("M61", -1085.4615)
    ,("H74", 556163.5600000001)
    ,("I74", 857998)
    ,("K74", 857998)
    ,("M74", -857124)
    ,("H77", 556163.5600000001)
    ,("I77", 857998)
    ,("K77", 4289990)
    ,("M77", -4285622)
    ,("H116", 2696.4443)
    ,("I116", 2769.261)
    ,("J116", 2277.75)
    ,("K116", 2769.261)
    ,("L116", 2277.75)
    ,("M116", -475.261)
    ,("N116", -6865.75)
    ,("H122", 1956.5)
    ,("I122", 1668.3334)
    ,("K122", 5005.0002)
    ,("M122", -2555.0002)
    ,("H132", 1758.1428)
    ,("I132", 1123.7778)
    ,("K132", 3371.3334)
    ,("M132", -841.3334000000004)
    ,("H136", 7285.2354)
    ,("I136", 1297.4615)
    ,("K136", 3892.3845)
    ,("M136", -1342.3845)
)
foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    ,("H3", 2696.4443)
    ,("I3", 2769.261)
    ,("J3", 2277.75)
    ,("K3", 2769.261)
    ,("L3", 2277.75)
    ,("M3", -2655.261)
    ,("N3", -2505.75)
    ,("H22", 999.6667)
    ,("J22", 0)
    ,("L22", 0)
    ,("H80", 789.5625)
    ,("I80", 921.5)
    ,("K80", 921.5)
    ,("M80", 76.5)
    ,("H83", 789.5625)
    ,("I83", 921.5)
    ,("K83", 4607.5)
    ,("M83", 384.5)
    ,("H99", 1078)
    ,("I99", 1078)
    ,("J99", 0)
    ,("K99", 1078)
    ,("L99", 0)
    ,("M99", 420)
)
foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
$clears = @("N22", "N99")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    ,("H5", 265.2)
    ,("I5", 42)
    ,("J5", 600)
    ,("K5", 42)
    ,("L5", 600)
    ,("M5", 70)
    ,("N5", -824)
    ,("H16", 2521.32)
    ,("I16", 2402.6316)
    ,("J16", 2897.1667)
    ,("K16", 2402.6316)
    ,("L16", 2897.1667)
    ,("M16", -2115.6316)
    ,("N16", -3471.1667)
    ,("H33", 2497.5)
    ,("I33", 2497.5)
    ,("K33", 2497.5)
    ,("M33", -2118.5)
    ,("H36", 9274.5)
    ,("I36", 9274.5)
    ,("K36", 9274.5)
    ,("M36", -8886.5)
    ,("H40", 9274.5)
    ,("I40", 9274.5)
    ,("K40", 9274.5)
    ,("M40", -9114.5)
    ,("H54", 44000)
    ,("J54", 44000)
    ,("L54", 44000)
    ,("N54", -45316)
    ,("H58", 1322.5264)
    ,("I58", 1075.8572)
    ,("J58", 2013.2)
    ,("K58", 1075.8572)
    ,("L58", 2013.2)
    ,("M58", -872.8571999999999)
    ,("N58", -2419.2)
    ,("H94", 2483.4211)
    ,("J94", 3904.875)
    ,("L94", 3904.875)
    ,("N94", -4806.875)
    ,("H113", 2521.32)
    ,("I113", 2402.6316)
    ,("J113", 2897.1667)
    ,("K113", 2402.6316)
    ,("L113", 2897.1667)
    ,("M113", -232.6316000000002)
    ,("N113", -7237.1667)
    ,("H132", 145429.28)
    ,("I132", 250752.5)
    ,("K132", 752257.5)
    ,("M132", -749727.5)
    ,("H136", 1322.5264)
    ,("I136", 1075.8572)
    ,("J136", 2013.2)
    ,("K136", 3227.5716)
    ,("L136", 6039.6)
    ,("M136", -677.5715999999998)
    ,("N136", -11139.6)
)
foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    ,("H12", 233.57895)
    ,("J12", 247.46153)
    ,("L12", 742.38459)
    ,("N12", -1088.38459)
    ,("H131", 101020.58)
    ,("J131", 1653.4193)
    ,("L131", 4960.257900000001)
    ,("N131", -15040.2579)
    ,("H134", 0)
    ,("I134", 0)
    ,("K134", 0)
)
foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
$clears = @("M134")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    ,("H41", 10474)
    ,("I41", 1998)
    ,("K41", 1998)
    ,("M41", -1643)
    ,("H113", 3044.8572)
    ,("I113", 2846.4375)
    ,("K113", 2846.4375)
    ,("M113", -676.4375)
    ,("H122", 4160.6523)
    ,("I122", 3949.7334)
    ,("J122", 4556.125)
    ,("K122", 11849.2002)
    ,("L122", 13668.375)
    ,("M122", -9399.200199999999)
    ,("N122", -18568.375)
    ,("H132", 1810.5385)
    ,("I132", 1565.75)
    ,("J132", 2202.2)
    ,("K132", 4697.25)
    ,("L132", 6606.599999999999)
    ,("M132", -2167.25)
    ,("N132", -11666.6)
)
foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    ,("H52", 48999.668)
    ,("I52", 0)
    ,("J52", 48999.668)
    ,("K52", 0)
    ,("L52", 48999.668)
    ,("N52", -49451.668)
    ,("H58", 21994.334)
    ,("I58", 11491.5)
    ,("K58", 11491.5)
    ,("M58", -11183.5)
    ,("H81", 6664.375)
    ,("I81", 7045)
    ,("J81", 4000)
    ,("K81", 14090)
    ,("L81", 8000)
    ,("M81", -13029)
    ,("N81", -10122)
    ,("H84", 6664.375)
    ,("I84", 7045)
    ,("J84", 4000)
    ,("K84", 70450)
    ,("L84", 40000)
    ,("M84", -65146)
    ,("N84", -50608)
    ,("H122", 10029159)
    ,("I122", 11938942)
    ,("J122", 2801.25)
    ,("K122", 35816826)
    ,("L122", 8403.75)
    ,("M122", -35814376)
    ,("N122", -13303.75)
    ,("H132", 5831791)
    ,("I132", 6598399)
    ,("J132", 5570)
    ,("K132", 19795197)
    ,("L132", 16710)
    ,("M132", -19792667)
    ,("N132", -21770)
)
foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
$clears = @("M52")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}
